# Update the QuantitativeMetrics evaluation sheet with new test-run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# --- Execution metrics block ---

# Compilation success: yes -> no
$ws.Range("B5").Value = "no"

# Runtime without error: yes -> (cleared)
$ws.Range("B6").ClearContents()

# Assertion validity: yes -> (cleared), and drop its note
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# --- Syntax & Semantic similarity block ---

# Code BLEU score updated
$ws.Range("B12").Value = 0.2785910815676469
$ws.Range("C12").Value = "{'codebleu': 0.2785910815676469, 'ngram_match_score': 0.06349531214452699, 'weighted_ngram_match_score': 0.1067298199868665, 'syntax_match_score': 0.5989010989010989, 'dataflow_match_score': 0.34523809523809523}"

# Move the active selection to B6, matching the updated review focus.
$ws.Activate()
$ws.Range("B6").Select()
